$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the placeholder empty inline-string cells left over from the
# template (row 2 B:G, plus a few scattered blanks in D/E/B/C across
# rows 3-7) now that real data has been filled in.
$ws.Range("B2:G2").ClearContents()
$ws.Range("D3:E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("B5:C5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()

# Fill in the "% Giam gia" (discount) column with the computed discount
# rates, formatted with a dark-red font to flag the decrease.
$ws.Range("D4").Value = -0.2320185614849146
$ws.Range("D5").Value = -0.08936550491509543
$ws.Range("D6").Value = -0.540540540540535
$ws.Range("D7").Value = -0.4739336492890933

$ws.Range("D4:D7").Font.Color = 170
